$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect to write, then restore protection afterwards.
$ws.Unprotect()

# Update H2 from "Market Mode/Surya Mitra/NBCFDC" to "Any"
$ws.Range("H2").Value = "Any"

# Move the active selection to H2
$ws.Range("H2").Select()

# Restore sheet protection
$ws.Protect()
